$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 30546.666
$ws.Range("I69").Value = 6044
$ws.Range("J69").Value = 48048.57
$ws.Range("K69").Value = 18132
$ws.Range("L69").Value = 144145.71
$ws.Range("M69").Value = -17258
$ws.Range("N69").Value = -145893.71

$ws.Range("H72").Value = 30546.666
$ws.Range("I72").Value = 6044
$ws.Range("J72").Value = 48048.57
$ws.Range("K72").Value = 54396
$ws.Range("L72").Value = 432437.13
$ws.Range("M72").Value = -50028
$ws.Range("N72").Value = -441173.13

$ws.Range("H92").Value = 1047.7391
$ws.Range("I92").Value = 985.4666999999999
$ws.Range("K92").Value = 985.4666999999999
$ws.Range("M92").Value = 262.5333000000001

$ws.Range("H107").Value = 185.0625
$ws.Range("I107").Value = 185.0625
$ws.Range("K107").Value = 185.0625
$ws.Range("M107").Value = 1734.9375

$ws.Range("H111").Value = 6209.091
$ws.Range("J111").Value = 8254.4
$ws.Range("L111").Value = 24763.2
$ws.Range("N111").Value = -30897.2

$ws.Range("H129").Value = 495836.34
$ws.Range("I129").Value = 557359.6
$ws.Range("J129").Value = 3650
$ws.Range("K129").Value = 1672078.8
$ws.Range("L129").Value = 10950
$ws.Range("M129").Value = -1667078.8
$ws.Range("N129").Value = -20950

$ws.Range("H132").Value = 24888.818
$ws.Range("I132").Value = 1570.2667
$ws.Range("J132").Value = 74857.14
$ws.Range("K132").Value = 4710.800099999999
$ws.Range("L132").Value = 224571.42
$ws.Range("M132").Value = -2180.800099999999
$ws.Range("N132").Value = -229631.42

$ws.Range("H138").Value = 2352.2534
$ws.Range("I138").Value = 1481.4
$ws.Range("J138").Value = 2693.7646
$ws.Range("K138").Value = 4444.200000000001
$ws.Range("L138").Value = 8081.293799999999
$ws.Range("M138").Value = 695.7999999999993
$ws.Range("N138").Value = -18361.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8704.242
$ws.Range("I32").Value = 7098.0703
$ws.Range("K32").Value = 7098.0703
$ws.Range("M32").Value = -6811.0703

$ws.Range("H45").Value = 2437.625
$ws.Range("I45").Value = 1873
$ws.Range("J45").Value = 2876.7778
$ws.Range("K45").Value = 1873
$ws.Range("L45").Value = 2876.7778
$ws.Range("M45").Value = -1496
$ws.Range("N45").Value = -3630.7778

$ws.Range("H63").Value = 2031.6
$ws.Range("I63").Value = 1882.6923
$ws.Range("K63").Value = 1882.6923
$ws.Range("M63").Value = -1196.6923

$ws.Range("H66").Value = 2031.6
$ws.Range("I66").Value = 1882.6923
$ws.Range("K66").Value = 9413.461499999999
$ws.Range("M66").Value = -5981.461499999999

$ws.Range("H97").Value = 1846.2858
$ws.Range("J97").Value = 2059.6
$ws.Range("L97").Value = 2059.6
$ws.Range("N97").Value = -3051.6

$ws.Range("H122").Value = 5243.0454
$ws.Range("I122").Value = 3280.5386
$ws.Range("K122").Value = 9841.6158
$ws.Range("M122").Value = -7391.6158

$ws.Range("H132").Value = 5429.769
$ws.Range("I132").Value = 5063.6665
$ws.Range("K132").Value = 15190.9995
$ws.Range("M132").Value = -12660.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 49999.953
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 49999.953
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 49999.953
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -52495.953

$ws.Range("H88").Value = 41666.332
$ws.Range("J88").Value = 41666.332
$ws.Range("L88").Value = 41666.332
$ws.Range("N88").Value = -42478.332

$ws.Range("H90").Value = 49999.953
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 49999.953
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 149999.859
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -162479.859

$ws.Range("H91").Value = 41666.332
$ws.Range("J91").Value = 41666.332
$ws.Range("L91").Value = 41666.332
$ws.Range("N91").Value = -44474.332

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 3389.7896
$ws.Range("I134").Value = 2525.375
$ws.Range("K134").Value = 7576.125
$ws.Range("M134").Value = -5041.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6747.9443
$ws.Range("J94").Value = 1500.7273
$ws.Range("L94").Value = 1500.7273
$ws.Range("N94").Value = -2402.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 36299.89
$ws.Range("I59").Value = 250
$ws.Range("J59").Value = 40806.125
$ws.Range("K59").Value = 750
$ws.Range("L59").Value = 122418.375
$ws.Range("M59").Value = -210
$ws.Range("N59").Value = -123498.375

$ws.Range("H122").Value = 1491.25
$ws.Range("J122").Value = 1596.0834
$ws.Range("L122").Value = 14364.7506
$ws.Range("N122").Value = -19264.7506

$ws.Range("H129").Value = 1940
$ws.Range("J129").Value = 2397.5
$ws.Range("L129").Value = 7192.5
$ws.Range("N129").Value = -17192.5

$ws.Range("H132").Value = 1217.1
$ws.Range("I132").Value = 1120.75
$ws.Range("J132").Value = 1361.625
$ws.Range("K132").Value = 10086.75
$ws.Range("L132").Value = 12254.625
$ws.Range("M132").Value = -7556.75
$ws.Range("N132").Value = -17314.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 340.44446
$ws.Range("I97").Value = 309.5
$ws.Range("J97").Value = 365.2
$ws.Range("K97").Value = 309.5
$ws.Range("L97").Value = 365.2
$ws.Range("M97").Value = 186.5
$ws.Range("N97").Value = -1357.2

$ws.Range("H102").Value = 1404.3334
$ws.Range("I102").Value = 1239.8334
$ws.Range("K102").Value = 1239.8334
$ws.Range("M102").Value = 382.1666

$ws.Range("H113").Value = 6154.7334
$ws.Range("I113").Value = 3589.6
$ws.Range("J113").Value = 11285
$ws.Range("K113").Value = 3589.6
$ws.Range("L113").Value = 11285
$ws.Range("M113").Value = -1419.6
$ws.Range("N113").Value = -15625

$ws.Range("H126").Value = 2726.2666
$ws.Range("I126").Value = 1914.9231
$ws.Range("K126").Value = 5744.7693
$ws.Range("M126").Value = -3274.7693

$ws.Range("H132").Value = 4330.75
$ws.Range("I132").Value = 2700.8572
$ws.Range("K132").Value = 8102.571599999999
$ws.Range("M132").Value = -5572.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 500.14285
$ws.Range("I16").Value = 512.6
$ws.Range("J16").Value = 469
$ws.Range("K16").Value = 512.6
$ws.Range("L16").Value = 469
$ws.Range("M16").Value = -342.6
$ws.Range("N16").Value = -809

$ws.Range("H93").Value = 911258.25
$ws.Range("I93").Value = 2355.25
$ws.Range("K93").Value = 2355.25
$ws.Range("M93").Value = -1107.25

$ws.Range("H122").Value = 5501.3335
$ws.Range("I122").Value = 2804
$ws.Range("K122").Value = 8412
$ws.Range("M122").Value = -5962

$ws.Range("H136").Value = 5049.1304
$ws.Range("I136").Value = 3253.2144
$ws.Range("K136").Value = 9759.643199999999
$ws.Range("M136").Value = -7209.643199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 12496
$ws.Range("J76").Value = 12496
$ws.Range("L76").Value = 12496
$ws.Range("N76").Value = -13126

$ws.Range("H79").Value = 12496
$ws.Range("J79").Value = 12496
$ws.Range("L79").Value = 12496
$ws.Range("N79").Value = -14680

$ws.Range("H81").Value = 7905.8125
$ws.Range("I81").Value = 18207.5
$ws.Range("K81").Value = 36415
$ws.Range("M81").Value = -35354

$ws.Range("H84").Value = 7905.8125
$ws.Range("I84").Value = 18207.5
$ws.Range("K84").Value = 182075
$ws.Range("M84").Value = -176771

$ws.Range("H136").Value = 2963.9697
$ws.Range("I136").Value = 1776.6957
$ws.Range("K136").Value = 5330.0871
$ws.Range("M136").Value = -2780.0871
